## trim the Request url
## Applies targeted edits to the catalog.xlsx workbook:
##  - Adds two new jsonAttributes payload strings (used by L3/L4 on Sheet1)
##  - Adds a new wrap-text style and applies it to a newly-created empty cell I4 on Sheet1
##  - Adjusts row 3 height on Sheet1
##  - Updates the saved view/selection state on both worksheets

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item("Sheet1")
$backup = $wb.Worksheets.Item("Backup")

# --- Sheet1: update jsonAttributes (column L) values for rows 3 and 4 ---
$sheet1.Range("L3").Value = "[`n    {`"actorId`": `"`${actorId}`"},`n   { `"actorName`":`"`${actorName}`"},`n   { `"actorType`": `"Vendor`"},`n   `"catalogId`",`n   { `"products`": []},`n   { `"catalogs`":null}`n]"

$sheet1.Range("L4").Value = "[`n    {`"actorId`": `"`${actorId}`"},`n    {`"actorName`": `"`${actorName}`"},`n    {`"actorType`": `"Vendor`"},`n    `"catalogId`",`n    {`"products`": []},`n   { `"catalogs`": null}`n]"

# --- Sheet1: new empty, word-wrapped cell at I4 ---
$sheet1.Range("I4").Value = $null
$sheet1.Range("I4").WrapText = $true

# --- Sheet1: row 3 height changes from 240 to 120 ---
$sheet1.Rows.Item(3).RowHeight = 120

# --- Sheet1: view/selection state ---
$sheet1.Activate()
$excel.ActiveWindow.Zoom = 100
$sheet1.Range("B4").Select()

# --- Backup: view/selection state (scroll so C3 is the top-left visible
#     cell, then select M3) ---
$backup.Activate()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 3
$backup.Range("M3").Select()

# --- Re-activate Sheet1 so it remains the selected/visible tab on save ---
$sheet1.Activate()
